$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("I2:I27")
$range.Value = $false
$range.NumberFormat = '"TRUE";"TRUE";"FALSE"'
$range.Select() | Out-Null
